# daily auto push: 2026-02-14 18:51 UTC
#
# A new reading was recorded for 2026/02/15 (日, hour 0, ranking 201).
# It belongs right before the existing 2026/12/29 block, so insert a new
# row at row 821 (shifting all following rows down by one) and populate it.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift row 821 (and everything below it) down by one row.
$ws.Rows("821:821").Insert()

# Column A holds dates as plain text (e.g. "2026/12/29"), not real Excel
# dates. Assigning a "yyyy/mm/dd"-shaped string straight to .Value makes
# Excel's input parser auto-convert it into a date serial (and tag the
# cell with a number-format style), which the source workbook does not
# use anywhere in this column. To keep the cell a literal text value with
# no extra style, stage the text (with a trailing space so the date
# parser leaves it alone), trim it with a formula, then paste the trimmed
# result back as a value.
$ws.Cells.Item(900, 1).Value = "2026/02/15 "
$ws.Cells.Item(900, 2).Formula = "=TRIM(A900)"
$ws.Cells.Item(900, 2).Copy()
$ws.Cells.Item(821, 1).PasteSpecial(-4163)  # xlPasteValues

# Remaining columns are plain text / numbers and don't trigger any
# special parsing.
$ws.Cells.Item(821, 2).Value = "日"
$ws.Cells.Item(821, 3).Value = 0
$ws.Cells.Item(821, 4).Value = 201

# Clean up the scratch cells so they don't linger in the saved sheet.
$ws.Cells.Item(900, 1).Value = ""
$ws.Cells.Item(900, 2).Value = ""
$excel.CutCopyMode = $false
